# Weekly update: insert the newest Jengibre price record at the top of the
# data block (row 134), pushing the existing rows down by one. All other
# rows keep their values; only the sheet grows by one row (R152 -> R153).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 134:152 down to 135:153 and open up a blank row 134.
$ws.Rows("134:134").Insert()

# Fill the new row 134 with this week's record.
$ws.Range("A134").Value = 8
$ws.Range("B134").Value = "Terminal La Palmera de La Serena"
$ws.Range("C134").Value = "Coquimbo"
$ws.Range("D134").Value = 45124
$ws.Range("E134").Value = 4
$ws.Range("F134").Value = 100114007
$ws.Range("G134").Value = "Jengibre"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 380
$ws.Range("K134").Value = 17500
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = 17750
$ws.Range("N134").Value = '$/caja 13 kilos'
$ws.Range("O134").Value = "Perú"
$ws.Range("P134").Value = 1365
$ws.Range("Q134").Value = 13
$ws.Range("R134").Value = "Hortaliza"
